# Applies the data correction described in the commit message:
#  - Column A (date/id) for rows 3..63 is shifted from 2015xxxx to 2017xxxx (+20000)
#  - Column E (value) for most of the same rows is corrected to a new value
#    (a handful of rows keep their original E value, matching the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column E, keyed by row number. Rows not present here keep
# their original E value (rows 9, 12, 26, 30, 41, 57).
$eValues = @{
    3=8; 4=9; 5=17; 6=19; 7=16; 8=19; 10=15; 11=15; 13=14; 14=14; 15=7;
    16=13; 17=15; 18=17; 19=9; 20=9; 21=12; 22=5; 23=15; 24=10; 25=7;
    27=14; 28=15; 29=12; 31=11; 32=10; 33=18; 34=15; 35=19; 36=7; 37=12;
    38=20; 39=12; 40=19; 42=8; 43=15; 44=6; 45=9; 46=6; 47=11; 48=10;
    49=10; 50=8; 51=10; 52=5; 53=11; 54=15; 55=17; 56=6; 58=17; 59=9;
    60=5; 61=20; 62=19; 63=8
}

for ($row = 3; $row -le 63; $row++) {
    $aCell = $ws.Cells.Item($row, 1)
    $aOld = $aCell.Value2
    if ($aOld -ne $null) {
        $aCell.Value2 = $aOld + 20000
    }

    if ($eValues.ContainsKey($row)) {
        $ws.Cells.Item($row, 5).Value2 = $eValues[$row]
    }
}
